# Corrections to the E2PROM Container Review Checklist "Checklist" sheet.
# Rows 18-24 are renumbered / cleaned up:
#   - old TestData_08..TestData_11 (rows 18-21) shift down to TestData_07..TestData_10
#   - a brand-new, cleaned-up TestData_11 row is inserted at row 22
#   - DUMMY_TestModuleCnt / ASDFClockTower move down one row (23/24) unchanged
#   - the old bottom "TestData_07" junk-data row (formerly row 24) is removed
#   - every row's scratch/rating columns (C,D,E,F,N,P,Q,R,S, stray G/H/I) are
#     wiped back to the single clean "Reprog" checkbox + description text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# Clear the full data block first (A:S, rows 18-24) so stray old values
# (ratings, reviewer initials, extra use-case marks, etc.) don't linger.
$ws.Range("A18:S24").ClearContents()

# Row 18: TestData_07
$ws.Range("A18").Value = "TestData_07"
$ws.Range("B18").Value = 12345
$ws.Range("G18").Value = "X"
$ws.Range("H18").Value = "X"
$ws.Range("M18").Value = "ee_range"
$ws.Range("O18").Value = "description=- Component: Test`n- REPROG info: To be evaluated."

# Row 19: TestData_08
$ws.Range("A19").Value = "TestData_08"
$ws.Range("B19").Value = 12346
$ws.Range("I19").Value = "X"
$ws.Range("M19").Value = "ee_range"
$ws.Range("O19").Value = "description=- Component: TST Data`n- REPROG info: undefined"

# Row 20: TestData_09
$ws.Range("A20").Value = "TestData_09"
$ws.Range("B20").Value = 12347
$ws.Range("G20").Value = "X"
$ws.Range("M20").Value = "ee_range"
$ws.Range("O20").Value = "description=- Component: TST Data`n- REPROG info: tbd"

# Row 21: TestData_10
$ws.Range("A21").Value = "TestData_10"
$ws.Range("B21").Value = 12348
$ws.Range("H21").Value = "X"
$ws.Range("I21").Value = "X"
$ws.Range("M21").Value = "ee_range"
$ws.Range("O21").Value = "description=- Component: TST`n- REPROG info: t.b.d"

# Row 22: TestData_11 (new row)
$ws.Range("A22").Value = "TestData_11"
$ws.Range("B22").Value = 12349
$ws.Range("G22").Value = "X"
$ws.Range("H22").Value = "X"
$ws.Range("I22").Value = "X"
$ws.Range("M22").Value = "ee_range"
$ws.Range("O22").Value = "description=- Component: TST`n- REPROG info: use case REPROG must be set"

# Row 23: DUMMY_TestModuleCnt (shifted down from row 22)
$ws.Range("A23").Value = "DUMMY_TestModuleCnt"
$ws.Range("B23").Value = 31416
$ws.Range("I23").Value = "X"
$ws.Range("M23").Value = "ee_erase"
$ws.Range("O23").Value = "description=- Component: DUMMY`n- REPROG info: use case REPROG must not be set (data must not be changed after reprogramming)!`nSometimes the description is longer than two rows.`nOther times, there are more than three.`nIn this case, it is one more than four. And could be more."

# Row 24: ASDFClockTower (shifted down from row 23)
$ws.Range("A24").Value = "ASDFClockTower"
$ws.Range("B24").Value = 111255
$ws.Range("G24").Value = "X"
$ws.Range("M24").Value = "ee_datablock"
$ws.Range("O24").Value = "description=- Component: ASDF`n- REPROG info: use case REPROG must not be set.`n- REPROG info: In certain cases there are two comments of this type.`nThere are also strings up to 160 characters per row, only on description fields and usually is not only one row. Like this example but a little bit longer."
